$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Memory Map")

# Insert 6 new rows before row 43 (shifting old row 43 "not used - expansion" and
# everything below it down by 6 rows), mirroring Excel's own row-insert behaviour
# so that all the dependent formulas / shared-formula ranges below get adjusted
# automatically.
$ws.Range("A43:N48").EntireRow.Insert()

# Carry the row-42 (VGAcols) formatting down onto the 6 freshly inserted blank
# rows so the new register rows look like the rest of the table (DEC column
# style, the "Name"/"R/W" columns, etc.)
$ws.Range("A42:I42").Copy()
$ws.Range("A43:I48").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the six new rows (43-48) describing the newly mapped MAC controller
# registers, following the same pattern used by the surrounding rows.
$macRows = @(
    @{ Row = 43; Name = "MACreadyRX" },
    @{ Row = 44; Name = "MACdataRX" },
    @{ Row = 45; Name = "MACchecksum_error" },
    @{ Row = 46; Name = "MACreadyTX" },
    @{ Row = 47; Name = "MACdataTX" },
    @{ Row = 48; Name = "MACtransmit_request" }
)

foreach ($item in $macRows) {
    $r = $item.Row
    $prev = $r - 1

    $ws.Cells.Item($r, 1).Value = 1
    $ws.Cells.Item($r, 2).Formula = "=B$prev+4"
    $ws.Cells.Item($r, 3).Formula = "=B$r+A$r-1"
    $ws.Cells.Item($r, 4).Formula = "=DEC2HEX(B$r,6)"
    $ws.Cells.Item($r, 5).Formula = "=DEC2HEX(C$r,6)"
    $ws.Cells.Item($r, 7).Value = $item.Name
    $ws.Cells.Item($r, 8).Value = "R/W"
    $ws.Cells.Item($r, 9).Value = $ws.Cells.Item(42, 9).Value()
}

# The row that used to be "43" (the "not used - expansion" bucket) is now row
# 49; its starting address should continue on from the last new MAC register
# (row 48) rather than the old neighbour (row 42).
$ws.Cells.Item(49, 2).Formula = "=B48+4"

# Restore the view state recorded in the author's commit (scrolled down to the
# newly added rows, with B50 selected).
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("B50").Select()

$wb.Save()
